# Apply "first round of tests" update to the Testing Round #1 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Round #1 ")

# --- Third tester (Colton) score columns (M/N) for each question block ---
# Block 1 (rows 5-9)
$ws.Range("M5").Value = 5
$ws.Range("N5").Value = 5
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 5
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 5
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 5
$ws.Range("M9").Value = 5
$ws.Range("N9").Value = 5

# Block 2 (rows 12-14)
$ws.Range("M12").Value = 5
$ws.Range("N12").Value = 5
$ws.Range("M13").Value = 5
$ws.Range("N13").Value = 3
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 4

# Block 3 (rows 17-21)
$ws.Range("M17").Value = 5
$ws.Range("N17").Value = 5
$ws.Range("M18").Value = 3
$ws.Range("N18").Value = 4
$ws.Range("M19").Value = 3
$ws.Range("N19").Value = 4
$ws.Range("M20").Value = 3
$ws.Range("N20").Value = 4
$ws.Range("M21").Value = 4
$ws.Range("N21").Value = 5

# --- New notes cell for the third tester, merged M24:P24 ---
$ws.Range("M24:P24").Merge()
$ws.Range("M24").Value = "The initial wave of testing went well, although limitations of feedback from the hardware of AdobeXD make some navigation difficult, namely due to the lacking of a scroll bar or feedback of hovering over buttons. Once the users got a feel for how the buttons worked and learned that they were a drop down menu, navigating the webpage to view the competitve streams was much more efficient. They also didnt find the home page button of <U.C> to be very intutitive."

$wb.Save()
